$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.209.33"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "2.897.51"
$ws.Range("E3").Value = "  +3.79%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "3.352.78"
$ws.Range("E15").Value = "  +3.72%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.919.72"
$ws.Range("E16").Value = "  +4.58%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E17").Value = "  +5.52%  "
$ws.Range("D18").Value = "52.212.27"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  +3.18%  "
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("E26").Value = "  +8.89%  "
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("E30").Value = "  +14.76%  "
$ws.Range("E31").Value = "  +9.27%  "
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("E33").Value = "  +11.74%  "
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("E36").Value = "  -11.83%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("E41").Value = "  +10.38%  "
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("E43").Value = "  +4.87%  "
$ws.Range("E44").Value = "  +6.22%  "
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").Value = "2.173.35"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("E49").Value = "  +22.60%  "
$ws.Range("E50").Value = "  +11.17%  "
$ws.Range("E51").Value = "  -0.90%  "

# Numeric-looking price strings must be forced to remain text.
# Use the classic apostrophe text-prefix, then clear the resulting
# quote-prefix formatting so the cell keeps its original (default) style.
$ws.Range("D5").Value = "'352.50"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'112.89"
$ws.Range("D6").ClearFormats()
$ws.Range("D10").Value = "'39.93"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.136"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.0861"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'19.82"
$ws.Range("D13").ClearFormats()
$ws.Range("D17").Value = "'0.993"
$ws.Range("D17").ClearFormats()
$ws.Range("D19").Value = "'7.62"
$ws.Range("D19").ClearFormats()
$ws.Range("D21").Value = "'14.13"
$ws.Range("D21").ClearFormats()
$ws.Range("D23").Value = "'70.83"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").Value = "'268.58"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").Value = "'2.78"
$ws.Range("D25").ClearFormats()
$ws.Range("D27").Value = "'26.71"
$ws.Range("D27").ClearFormats()
$ws.Range("D29").Value = "'10.58"
$ws.Range("D29").ClearFormats()
$ws.Range("D31").Value = "'6.66"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").Value = "'37.44"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").Value = "'6.18"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").Value = "'53.14"
$ws.Range("D34").ClearFormats()
$ws.Range("D36").Value = "'1.99"
$ws.Range("D36").ClearFormats()
$ws.Range("D38").Value = "'3.31"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").Value = "'18.83"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").Value = "'2.04"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'2.76"
$ws.Range("D41").ClearFormats()
$ws.Range("D43").Value = "'22.76"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "'2.59"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'119.81"
$ws.Range("D45").ClearFormats()
$ws.Range("D49").Value = "'0.264"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'0.0344"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").Value = "'0.948"
$ws.Range("D51").ClearFormats()
